# ------------------------------------------------------------------
# Dickson_N=3.xlsx - add the B = 0.05 measurement data
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Dickson_Real_N=3")
$ws2 = $wb.Worksheets.Item("Dickson_Model_N3")

# ====================================================================
# Sheet "Dickson_Real_N=3"  (sheet1)
# ====================================================================

# --- Row 22: new section label row (R = 1 / T = 0.05 / B = 0.05) ---
# Shared strings must be minted in this order: "T = 0.05", "B = 0.05"
# (before the later "T" / "B" column headers) to match the workbook's
# sharedStrings append order.
$ws1.Range("B22").Value = "T = 0.05"
$ws1.Range("C22").Value = "B = 0.05"

# --- Row 13: add the per-row T / B fraction columns to the existing
#     B = 0.01 table; the old "eta" header moves from F13 to H13 ---
$ws1.Range("H13").Value = $ws1.Range("F13").Value
$ws1.Range("F13").Value = "T"
$ws1.Range("G13").Value = "B"

# --- Rows 14-18: insert the per-row T/B fraction values (0.01) and
#     move the eta formula from column F to column H, using the
#     updated measured Iout/Iin values ---
$ws1.Range("B14").Value = 15.625
$ws1.Range("C14").Value = 67.95
$ws1.Range("F14").Value = 0.01
$ws1.Range("G14").Value = 0.01
$ws1.Range("H14").Formula = "=B14*D14/(C14*E14)"

$ws1.Range("B15").Value = 26.571000000000002
$ws1.Range("C15").Value = 117.17
$ws1.Range("F15").Value = 0.01
$ws1.Range("G15").Value = 0.01
$ws1.Range("H15").Formula = "=B15*D15/(C15*E15)"

$ws1.Range("B16").Value = 33.438000000000002
$ws1.Range("C16").Value = 155.5
$ws1.Range("F16").Value = 0.01
$ws1.Range("G16").Value = 0.01
$ws1.Range("H16").Formula = "=B16*D16/(C16*E16)"

$ws1.Range("B17").Value = 34.078000000000003
$ws1.Range("C17").Value = 168.89
$ws1.Range("F17").Value = 0.01
$ws1.Range("G17").Value = 0.01
$ws1.Range("H17").Formula = "=B17*D17/(C17*E17)"

$ws1.Range("B18").Value = 33.235999999999997
$ws1.Range("C18").Value = 176.37
$ws1.Range("F18").Value = 0.01
$ws1.Range("G18").Value = 0.01
$ws1.Range("H18").Formula = "=B18*D18/(C18*E18)"

# --- Rows 21-28: brand-new B = 0.05 block ---
$ws1.Range("A21").Value = "Cfly = 1n per core, 2 cores per stage"

$ws1.Range("A22").Value = "R = 1"
# B22/C22 set above

$ws1.Range("A23").Value = "f (MHz)"
$ws1.Range("B23").Value = "Iout (mA)"
$ws1.Range("C23").Value = "Iin (mA)"
$ws1.Range("D23").Value = "Vout"
$ws1.Range("E23").Value = "Vin"
$ws1.Range("F23").Value = "eta"

$ws1.Range("A24").Value = 50
$ws1.Range("B24").Value = 12.265000000000001
$ws1.Range("C24").Value = 76.5
$ws1.Range("D24").Value = 3.5
$ws1.Range("E24").Value = 1
$ws1.Range("F24").Formula = "=B24*D24/(C24*E24)"

$ws1.Range("A25").Value = 100
$ws1.Range("B25").Value = 19.567
$ws1.Range("C25").Value = 133.11000000000001
$ws1.Range("D25").Value = 3.5
$ws1.Range("E25").Value = 1

$ws1.Range("A26").Value = 200
$ws1.Range("B26").Value = 18.942
$ws1.Range("C26").Value = 185.43
$ws1.Range("D26").Value = 3.5
$ws1.Range("E26").Value = 1

$ws1.Range("A27").Value = 300
$ws1.Range("B27").Value = 12.17
$ws1.Range("C27").Value = 213.21
$ws1.Range("D27").Value = 3.5
$ws1.Range("E27").Value = 1

$ws1.Range("A28").Value = 400
$ws1.Range("B28").Value = 3.976
$ws1.Range("C28").Value = 235.2
$ws1.Range("D28").Value = 3.5
$ws1.Range("E28").Value = 1

# shared formula block F25:F28
$ws1.Range("F25:F28").Formula = "=B25*D25/(C25*E25)"

# ====================================================================
# Sheet "Dickson_Model_N3"  (sheet2)
# ====================================================================

# --- New first block (rows 9-16): B = 0.01 measurement ---
$ws2.Range("A9").Value = "Cfly = 1n per core, 2 cores per stage"

$ws2.Range("A10").Value = "R = 1"
$ws2.Range("B10").Value = "T = 0.01"
$ws2.Range("C10").Value = "B = 0.01"

$ws2.Range("A11").Value = "f (MHz)"
$ws2.Range("B11").Value = "Iout (mA)"
$ws2.Range("C11").Value = "Iin (mA)"
$ws2.Range("D11").Value = "Vout"
$ws2.Range("E11").Value = "Vin"
$ws2.Range("F11").Value = "eta"

$ws2.Range("A12").Value = 50
$ws2.Range("B12").Value = 15.598000000000001
$ws2.Range("C12").Value = 67.91
$ws2.Range("D12").Value = 3.5
$ws2.Range("E12").Value = 1
$ws2.Range("F12").Formula = "=B12*D12/(C12*E12)"

$ws2.Range("A13").Value = 100
$ws2.Range("B13").Value = 26.56
$ws2.Range("C13").Value = 117.035
$ws2.Range("D13").Value = 3.5
$ws2.Range("E13").Value = 1

$ws2.Range("A14").Value = 200
$ws2.Range("B14").Value = 33.340000000000003
$ws2.Range("C14").Value = 155.34
$ws2.Range("D14").Value = 3.5
$ws2.Range("E14").Value = 1

$ws2.Range("A15").Value = 300
$ws2.Range("B15").Value = 33.93
$ws2.Range("C15").Value = 168.77
$ws2.Range("D15").Value = 3.5
$ws2.Range("E15").Value = 1

$ws2.Range("A16").Value = 400
$ws2.Range("B16").Value = 33.04
$ws2.Range("C16").Value = 176.16
$ws2.Range("D16").Value = 3.5
$ws2.Range("E16").Value = 1

$ws2.Range("F13:F16").Formula = "=B13*D13/(C13*E13)"

# --- New second block (rows 18-25): B = 0.05 measurement ---
$ws2.Range("A18").Value = "Cfly = 1n per core, 2 cores per stage"

$ws2.Range("A19").Value = "R = 1"
$ws2.Range("B19").Value = "T = 0.05"
$ws2.Range("C19").Value = "B = 0.05"

$ws2.Range("A20").Value = "f (MHz)"
$ws2.Range("B20").Value = "Iout (mA)"
$ws2.Range("C20").Value = "Iin (mA)"
$ws2.Range("D20").Value = "Vout"
$ws2.Range("E20").Value = "Vin"
$ws2.Range("F20").Value = "eta"

$ws2.Range("A21").Value = 50
$ws2.Range("B21").Value = 12.24
$ws2.Range("C21").Value = 76.466999999999999
$ws2.Range("D21").Value = 3.5
$ws2.Range("E21").Value = 1
$ws2.Range("F21").Formula = "=B21*D21/(C21*E21)"

$ws2.Range("A22").Value = 100
$ws2.Range("B22").Value = 19.515999999999998
$ws2.Range("C22").Value = 133.03
$ws2.Range("D22").Value = 3.5
$ws2.Range("E22").Value = 1

$ws2.Range("A23").Value = 200
$ws2.Range("B23").Value = 18.838999999999999
$ws2.Range("C23").Value = 185.35
$ws2.Range("D23").Value = 3.5
$ws2.Range("E23").Value = 1

$ws2.Range("A24").Value = 300
$ws2.Range("B24").Value = 12.015000000000001
$ws2.Range("C24").Value = 213.05
$ws2.Range("D24").Value = 3.5
$ws2.Range("E24").Value = 1

$ws2.Range("A25").Value = 400
$ws2.Range("B25").Value = 3.7549999999999999
$ws2.Range("C25").Value = 235.02
$ws2.Range("D25").Value = 3.5
$ws2.Range("E25").Value = 1

$ws2.Range("F22:F25").Formula = "=B22*D22/(C22*E22)"

# ====================================================================
# View state: sheet1 ("Dickson_Real_N=3") ends up the active/selected
# tab with I16 selected; sheet2 scrolled with C22 selected.
# ====================================================================
$ws2.Select()
$ws2.Range("C22").Select()

$ws1.Select()
$ws1.Range("I16").Select()
